# Implement food vouchers for treated DS-TB and MDR-TB cases
#
# This script:
#  1. Appends 10 new parameter rows at the bottom of the "constants" sheet
#     for the food-voucher economics parameters (DS-TB group then MDR-TB
#     group), mirroring the existing "improve_dst" econ parameter block.
#  2. Inserts one new row near the top of the same data block
#     (program_prop_food_voucher_improvement) which pushes all the
#     subsequent rows down by one.
#  3. Tidies up a handful of cosmetic view/window settings to match the
#     saved state of the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# ---------------------------------------------------------------------
# Step 1: append the 10 new "food_voucher" econ parameter rows at the
# bottom of the sheet (rows 223-232 using the sheet's current, pre-insert
# numbering -- they become rows 224-233 once the row is inserted below).
# Setting these values FIRST ensures the new shared strings they
# introduce are allocated before the program_prop_food_voucher_improvement
# strings added in step 2, matching the target shared-string order.
# ---------------------------------------------------------------------

$lastRow = $ws.Range("A222")
$dummyText = $lastRow.Offset(0, 4).Value2

$newRows = @(
    @{ Name = "econ_unitcost_food_voucher_ds";         Value = 10 },
    @{ Name = "econ_inflectioncost_food_voucher_ds";    Value = 0 },
    @{ Name = "econ_startupcost_food_voucher_ds";       Value = 10 },
    @{ Name = "econ_startupduration_food_voucher_ds";   Value = 3 },
    @{ Name = "econ_saturation_food_voucher_ds";        Value = 1 },
    @{ Name = "econ_unitcost_food_voucher_mdr";         Value = 10 },
    @{ Name = "econ_inflectioncost_food_voucher_mdr";   Value = 0 },
    @{ Name = "econ_startupcost_food_voucher_mdr";      Value = 10 },
    @{ Name = "econ_startupduration_food_voucher_mdr";  Value = 3 },
    @{ Name = "econ_saturation_food_voucher_mdr";       Value = 1 }
)

$row = 223
foreach ($item in $newRows) {
    $ws.Cells.Item($row, 1).Value = $item.Name
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item(222, 1).Style
    $ws.Cells.Item($row, 2).Value = $item.Value
    $ws.Cells.Item($row, 2).Style = $ws.Cells.Item(222, 2).Style
    $ws.Cells.Item($row, 5).Value = "dummy value for the moment"
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Step 2: insert a new row at row 46 for program_prop_food_voucher_improvement,
# pushing the existing row 46 (and everything after it) down by one.
# ---------------------------------------------------------------------

$ws.Rows.Item(46).Insert()
$ws.Range("A46").Value = "program_prop_food_voucher_improvement"
$ws.Range("B46").Value = 0.2
$ws.Range("E46").Value = "Proportional reduction in adverse outcomes from the food voucher intervention"

# ---------------------------------------------------------------------
# Step 3: cosmetic view/window tweaks to mirror the saved workbook state.
# ---------------------------------------------------------------------

# Sheet1 ("constants"): zoom to 100%, move selection to E229.
$win = $wb.Windows.Item(1)
$win.Zoom = 100
$ws.Range("E229").Select()

# Sheet2 ("time_variants"): move frozen-pane anchor & selection.
$ws2 = $wb.Worksheets.Item("time_variants")
$ws2.Activate()
$ws2.Range("C2").Select()
$wb.Windows.Item(1).FreezePanes = $false
$wb.Windows.Item(1).FreezePanes = $true
$ws2.Range("A21").Select()

# Sheet3 ("dropdown_lists"): scroll so row 4 is the top-left visible cell.
$ws3 = $wb.Worksheets.Item("dropdown_lists")
$ws3.Activate()
$win.ScrollRow = 4

# Re-activate the constants sheet (tabSelected="1" in the target) and
# restore its selection/scroll position.
$ws.Activate()
$win.ScrollRow = 210
$ws.Range("E229").Select()

# Minimize the workbook window (workbookView minimized="1" in the target).
$excel.WindowState = -4140
